# Update column F ("dSF") values on the (only) worksheet to reflect a
# repull of data / recalculated mean, per the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    14 = -1
    17 = -2
    18 = 4
    19 = -4
    21 = 1
    23 = -7
    25 = -4
    26 = -1
    27 = 1
    30 = -9
    34 = -6
    36 = 5
    38 = 1
    39 = 5
    41 = -3
    45 = -6
    50 = -5
    52 = 0
    56 = -4
    61 = -2
    64 = 5
    66 = 5
    67 = -3
    69 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
